$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 3.9
$ws.Cells.Item(2, 7).Value = 4.3
$ws.Cells.Item(2, 8).Value = 2.28
$ws.Cells.Item(2, 9).Value = 2.44
$ws.Cells.Item(2, 10).Value = 2.94
$ws.Cells.Item(2, 11).Value = 3.15
$ws.Cells.Item(2, 12).Value = 1.61
$ws.Cells.Item(2, 13).Value = 1.13
$ws.Cells.Item(2, 15).Value = 1.57
$ws.Cells.Item(2, 16).Value = 1.54
$ws.Cells.Item(2, 17).Value = 2.68
$ws.Cells.Item(2, 18).Value = 1.19
$ws.Cells.Item(2, 19).Value = 5.6
$ws.Cells.Item(2, 20).Value = 2.1
$ws.Cells.Item(2, 21).Value = 1.82
$ws.Cells.Item(2, 22).Value = 1.7
$ws.Cells.Item(2, 23).Value = 1.3
$ws.Cells.Item(2, 24).Value = 8.6
$ws.Cells.Item(2, 25).Value = 7.6
$ws.Cells.Item(2, 26).Value = 14.5
$ws.Cells.Item(2, 27).Value = 40
$ws.Cells.Item(2, 28).Value = 11.5
$ws.Cells.Item(2, 29).Value = 6.8
$ws.Cells.Item(2, 30).Value = 12
$ws.Cells.Item(2, 31).Value = 36
$ws.Cells.Item(2, 32).Value = 28
$ws.Cells.Item(2, 33).Value = 18
$ws.Cells.Item(2, 34).Value = 22
$ws.Cells.Item(2, 35).Value = 85
$ws.Cells.Item(2, 36).Value = 110
$ws.Cells.Item(2, 37).Value = 90
$ws.Cells.Item(2, 38).Value = 240
$ws.Cells.Item(2, 40).Value = 600
$ws.Cells.Item(2, 41).Value = 36
$ws.Cells.Item(3, 6).Value = 5
$ws.Cells.Item(3, 7).Value = 5.3
$ws.Cells.Item(3, 8).Value = 1.77
$ws.Cells.Item(3, 9).Value = 1.81
$ws.Cells.Item(3, 10).Value = 3.9
$ws.Cells.Item(3, 11).Value = 4.1
$ws.Cells.Item(3, 12).Value = 1.47
$ws.Cells.Item(3, 13).Value = 1.09
$ws.Cells.Item(3, 14).Value = 3.5
$ws.Cells.Item(3, 15).Value = 1.37
$ws.Cells.Item(3, 16).Value = 1.84
$ws.Cells.Item(3, 17).Value = 2.12
$ws.Cells.Item(3, 18).Value = 1.3
$ws.Cells.Item(3, 19).Value = 4
$ws.Cells.Item(3, 20).Value = 2
$ws.Cells.Item(3, 21).Value = 1.86
$ws.Cells.Item(3, 22).Value = 2.2
$ws.Cells.Item(3, 23).Value = 1.23
$ws.Cells.Item(3, 24).Value = 13
$ws.Cells.Item(3, 25).Value = 7.6
$ws.Cells.Item(3, 26).Value = 10.5
$ws.Cells.Item(3, 27).Value = 19
$ws.Cells.Item(3, 28).Value = 16
$ws.Cells.Item(3, 29).Value = 9
$ws.Cells.Item(3, 30).Value = 9.8
$ws.Cells.Item(3, 31).Value = 21
$ws.Cells.Item(3, 32).Value = 60
$ws.Cells.Item(3, 33).Value = 21
$ws.Cells.Item(3, 34).Value = 22
$ws.Cells.Item(3, 35).Value = 46
$ws.Cells.Item(3, 36).Value = 140
$ws.Cells.Item(3, 37).Value = 80
$ws.Cells.Item(3, 38).Value = 90
$ws.Cells.Item(3, 39).Value = 580
$ws.Cells.Item(3, 40).Value = 200
$ws.Cells.Item(3, 41).Value = 14.5
$ws.Cells.Item(4, 6).Value = 1.38
$ws.Cells.Item(4, 7).Value = 1.39
$ws.Cells.Item(4, 8).Value = 11
$ws.Cells.Item(4, 9).Value = 11.5
$ws.Cells.Item(4, 10).Value = 5.2
$ws.Cells.Item(4, 11).Value = 5.3
$ws.Cells.Item(4, 12).Value = 1.37
$ws.Cells.Item(4, 13).Value = 1.06
$ws.Cells.Item(4, 14).Value = 4.3
$ws.Cells.Item(4, 16).Value = 2.12
$ws.Cells.Item(4, 17).Value = 1.86
$ws.Cells.Item(4, 18).Value = 1.43
$ws.Cells.Item(4, 19).Value = 3.2
$ws.Cells.Item(4, 20).Value = 2.28
$ws.Cells.Item(4, 21).Value = 1.75
$ws.Cells.Item(4, 22).Value = 1.09
$ws.Cells.Item(4, 23).Value = 3.55
$ws.Cells.Item(4, 24).Value = 17
$ws.Cells.Item(4, 25).Value = 29
$ws.Cells.Item(4, 26).Value = 100
$ws.Cells.Item(4, 27).Value = 470
$ws.Cells.Item(4, 29).Value = 11.5
$ws.Cells.Item(4, 30).Value = 38
$ws.Cells.Item(4, 31).Value = 220
$ws.Cells.Item(4, 32).Value = 7.2
$ws.Cells.Item(4, 34).Value = 32
$ws.Cells.Item(4, 35).Value = 200
$ws.Cells.Item(4, 36).Value = 11
$ws.Cells.Item(4, 38).Value = 48
$ws.Cells.Item(4, 39).Value = 210
$ws.Cells.Item(4, 40).Value = 6.2
$ws.Cells.Item(4, 41).Value = 290
$ws.Cells.Item(5, 6).Value = 1.8
$ws.Cells.Item(5, 7).Value = 1.85
$ws.Cells.Item(5, 8).Value = 5.6
$ws.Cells.Item(5, 9).Value = 6
$ws.Cells.Item(5, 10).Value = 3.5
$ws.Cells.Item(5, 14).Value = 2.92
$ws.Cells.Item(5, 15).Value = 1.49
$ws.Cells.Item(5, 16).Value = 1.63
$ws.Cells.Item(5, 17).Value = 2.52
$ws.Cells.Item(5, 18).Value = 1.22
$ws.Cells.Item(5, 19).Value = 5
$ws.Cells.Item(5, 21).Value = 1.71
$ws.Cells.Item(5, 22).Value = 1.2
$ws.Cells.Item(5, 23).Value = 2.16
$ws.Cells.Item(5, 25).Value = 15.5
$ws.Cells.Item(5, 26).Value = 42
$ws.Cells.Item(5, 27).Value = 180
$ws.Cells.Item(5, 28).Value = 6.6
$ws.Cells.Item(5, 29).Value = 8.4
$ws.Cells.Item(5, 30).Value = 24
$ws.Cells.Item(5, 31).Value = 110
$ws.Cells.Item(5, 32).Value = 9
$ws.Cells.Item(5, 34).Value = 980
$ws.Cells.Item(5, 35).Value = 200
$ws.Cells.Item(5, 36).Value = 19
$ws.Cells.Item(5, 37).Value = 24
$ws.Cells.Item(5, 39).Value = 230
$ws.Cells.Item(5, 40).Value = 18.5
$ws.Cells.Item(5, 41).Value = 180
$ws.Cells.Item(6, 6).Value = 2.6
$ws.Cells.Item(6, 8).Value = 2.82
$ws.Cells.Item(6, 9).Value = 3.05
$ws.Cells.Item(6, 11).Value = 3.6
$ws.Cells.Item(6, 14).Value = 3.6
$ws.Cells.Item(6, 15).Value = 1.33
$ws.Cells.Item(6, 16).Value = 1.92
$ws.Cells.Item(6, 19).Value = 3.65
$ws.Cells.Item(6, 20).Value = 1.74
$ws.Cells.Item(6, 21).Value = 2.06
$ws.Cells.Item(6, 23).Value = 1.55
$ws.Cells.Item(6, 24).Value = 90
$ws.Cells.Item(6, 25).Value = 21
$ws.Cells.Item(6, 27).Value = 900
$ws.Cells.Item(6, 30).Value = 23
$ws.Cells.Item(6, 31).Value = 1000
$ws.Cells.Item(6, 32).Value = 38
$ws.Cells.Item(6, 34).Value = 20
$ws.Cells.Item(6, 35).Value = 1000
$ws.Cells.Item(6, 36).Value = 130
$ws.Cells.Item(6, 37).Value = 55
$ws.Cells.Item(6, 39).Value = 580
$ws.Cells.Item(7, 6).Value = 1.94
$ws.Cells.Item(7, 8).Value = 4.6
$ws.Cells.Item(7, 10).Value = 3.35
$ws.Cells.Item(7, 11).Value = 3.6
$ws.Cells.Item(7, 14).Value = 3
$ws.Cells.Item(7, 15).Value = 1.45
$ws.Cells.Item(7, 16).Value = 1.66
$ws.Cells.Item(7, 17).Value = 2.36
$ws.Cells.Item(7, 18).Value = 1.24
$ws.Cells.Item(7, 20).Value = 2.08
$ws.Cells.Item(7, 23).Value = 1.98
$ws.Cells.Item(7, 25).Value = 23
$ws.Cells.Item(7, 26).Value = 55
$ws.Cells.Item(7, 27).Value = 1000
$ws.Cells.Item(7, 29).Value = 8.2
$ws.Cells.Item(7, 30).Value = 36
$ws.Cells.Item(7, 31).Value = 1000
$ws.Cells.Item(7, 35).Value = 1000
$ws.Cells.Item(7, 36).Value = 80
$ws.Cells.Item(7, 37).Value = 75
$ws.Cells.Item(7, 39).Value = 1000
$ws.Cells.Item(8, 6).Value = 3.85
$ws.Cells.Item(8, 7).Value = 3.9
$ws.Cells.Item(8, 8).Value = 2.36
$ws.Cells.Item(8, 9).Value = 2.38
$ws.Cells.Item(8, 10).Value = 3.1
$ws.Cells.Item(8, 11).Value = 3.2
$ws.Cells.Item(8, 12).Value = 1.63
$ws.Cells.Item(8, 13).Value = 1.15
$ws.Cells.Item(8, 14).Value = 2.52
$ws.Cells.Item(8, 15).Value = 1.64
$ws.Cells.Item(8, 16).Value = 1.5
$ws.Cells.Item(8, 17).Value = 2.92
$ws.Cells.Item(8, 18).Value = 1.17
$ws.Cells.Item(8, 19).Value = 6.4
$ws.Cells.Item(8, 20).Value = 2.44
$ws.Cells.Item(8, 21).Value = 1.66
$ws.Cells.Item(8, 22).Value = 1.72
$ws.Cells.Item(8, 23).Value = 1.34
$ws.Cells.Item(8, 24).Value = 7.4
$ws.Cells.Item(8, 25).Value = 7
$ws.Cells.Item(8, 26).Value = 12
$ws.Cells.Item(8, 27).Value = 29
$ws.Cells.Item(8, 28).Value = 9.6
$ws.Cells.Item(8, 29).Value = 7.6
$ws.Cells.Item(8, 30).Value = 13.5
$ws.Cells.Item(8, 31).Value = 36
$ws.Cells.Item(8, 32).Value = 23
$ws.Cells.Item(8, 34).Value = 26
$ws.Cells.Item(8, 35).Value = 75
$ws.Cells.Item(8, 36).Value = 85
$ws.Cells.Item(8, 38).Value = 100
$ws.Cells.Item(8, 39).Value = 220
$ws.Cells.Item(8, 40).Value = 130
$ws.Cells.Item(8, 41).Value = 36
$ws.Cells.Item(9, 6).Value = 1.98
$ws.Cells.Item(9, 7).Value = 1.99
$ws.Cells.Item(9, 8).Value = 4.5
$ws.Cells.Item(9, 9).Value = 4.6
$ws.Cells.Item(9, 10).Value = 3.65
$ws.Cells.Item(9, 12).Value = 1.46
$ws.Cells.Item(9, 13).Value = 1.08
$ws.Cells.Item(9, 14).Value = 3.6
$ws.Cells.Item(9, 15).Value = 1.36
$ws.Cells.Item(9, 16).Value = 1.88
$ws.Cells.Item(9, 17).Value = 2.08
$ws.Cells.Item(9, 18).Value = 1.33
$ws.Cells.Item(9, 19).Value = 3.85
$ws.Cells.Item(9, 20).Value = 1.92
$ws.Cells.Item(9, 21).Value = 2.04
$ws.Cells.Item(9, 22).Value = 1.28
$ws.Cells.Item(9, 23).Value = 2
$ws.Cells.Item(9, 24).Value = 13
$ws.Cells.Item(9, 25).Value = 15
$ws.Cells.Item(9, 26).Value = 32
$ws.Cells.Item(9, 27).Value = 95
$ws.Cells.Item(9, 28).Value = 8.4
$ws.Cells.Item(9, 29).Value = 8
$ws.Cells.Item(9, 30).Value = 17.5
$ws.Cells.Item(9, 31).Value = 60
$ws.Cells.Item(9, 32).Value = 11.5
$ws.Cells.Item(9, 33).Value = 10
$ws.Cells.Item(9, 34).Value = 18.5
$ws.Cells.Item(9, 35).Value = 65
$ws.Cells.Item(9, 36).Value = 22
$ws.Cells.Item(9, 37).Value = 21
$ws.Cells.Item(9, 38).Value = 38
$ws.Cells.Item(9, 39).Value = 120
$ws.Cells.Item(9, 40).Value = 18
$ws.Cells.Item(9, 41).Value = 65
$ws.Cells.Item(10, 6).Value = 2.08
$ws.Cells.Item(10, 7).Value = 2.1
$ws.Cells.Item(10, 8).Value = 3.8
$ws.Cells.Item(10, 9).Value = 3.9
$ws.Cells.Item(10, 10).Value = 3.8
$ws.Cells.Item(10, 11).Value = 3.85
$ws.Cells.Item(10, 12).Value = 1.37
$ws.Cells.Item(10, 13).Value = 1.06
$ws.Cells.Item(10, 14).Value = 4.4
$ws.Cells.Item(10, 15).Value = 1.27
$ws.Cells.Item(10, 16).Value = 2.16
$ws.Cells.Item(10, 17).Value = 1.83
$ws.Cells.Item(10, 18).Value = 1.45
$ws.Cells.Item(10, 19).Value = 3.15
$ws.Cells.Item(10, 20).Value = 1.72
$ws.Cells.Item(10, 21).Value = 2.32
$ws.Cells.Item(10, 22).Value = 1.34
$ws.Cells.Item(10, 23).Value = 1.91
$ws.Cells.Item(10, 24).Value = 17.5
$ws.Cells.Item(10, 25).Value = 16
$ws.Cells.Item(10, 26).Value = 28
$ws.Cells.Item(10, 27).Value = 75
$ws.Cells.Item(10, 28).Value = 11
$ws.Cells.Item(10, 29).Value = 8.4
$ws.Cells.Item(10, 32).Value = 13.5
$ws.Cells.Item(10, 34).Value = 16.5
$ws.Cells.Item(10, 35).Value = 48
$ws.Cells.Item(10, 36).Value = 24
$ws.Cells.Item(10, 37).Value = 20
$ws.Cells.Item(10, 38).Value = 34
$ws.Cells.Item(10, 39).Value = 80
$ws.Cells.Item(10, 40).Value = 13
$ws.Cells.Item(10, 41).Value = 36
$ws.Cells.Item(11, 6).Value = 2.68
$ws.Cells.Item(11, 7).Value = 2.78
$ws.Cells.Item(11, 8).Value = 2.86
$ws.Cells.Item(11, 9).Value = 3
$ws.Cells.Item(11, 11).Value = 3.45
$ws.Cells.Item(11, 14).Value = 4.3
$ws.Cells.Item(11, 15).Value = 1.28
$ws.Cells.Item(11, 16).Value = 2.06
$ws.Cells.Item(11, 18).Value = 1.43
$ws.Cells.Item(11, 19).Value = 3.15
$ws.Cells.Item(11, 20).Value = 1.62
$ws.Cells.Item(11, 21).Value = 2.36
$ws.Cells.Item(11, 22).Value = 1.5
$ws.Cells.Item(11, 23).Value = 1.56
$ws.Cells.Item(11, 24).Value = 15.5
$ws.Cells.Item(11, 25).Value = 13.5
$ws.Cells.Item(11, 26).Value = 19.5
$ws.Cells.Item(11, 27).Value = 46
$ws.Cells.Item(11, 28).Value = 13.5
$ws.Cells.Item(11, 29).Value = 8
$ws.Cells.Item(11, 30).Value = 12.5
$ws.Cells.Item(11, 32).Value = 19
$ws.Cells.Item(11, 33).Value = 12
$ws.Cells.Item(11, 34).Value = 16
$ws.Cells.Item(11, 35).Value = 40
$ws.Cells.Item(11, 36).Value = 65
$ws.Cells.Item(11, 37).Value = 29
$ws.Cells.Item(11, 38).Value = 48
$ws.Cells.Item(11, 39).Value = 80
$ws.Cells.Item(11, 40).Value = 21
$ws.Cells.Item(11, 41).Value = 24

Write-Output "Applied 307 cell updates"